$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2905.8
$ws.Range("I2").Value = 876.5
$ws.Range("J2").Value = 5949.75
$ws.Range("K2").Value = 876.5
$ws.Range("L2").Value = 5949.75
$ws.Range("M2").Value = -763.5
$ws.Range("N2").Value = -6175.75

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15496

$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16716

$ws.Range("H98").Value = 3086.0454
$ws.Range("J98").Value = 4852.4165
$ws.Range("L98").Value = 4852.4165
$ws.Range("N98").Value = -7848.4165

$ws.Range("H113").Value = 3012.5
$ws.Range("I113").Value = 2933.3333
$ws.Range("K113").Value = 2933.3333
$ws.Range("M113").Value = 320.6667000000002

$ws.Range("H122").Value = 3086.0454
$ws.Range("J122").Value = 4852.4165
$ws.Range("L122").Value = 14557.2495
$ws.Range("N122").Value = -19457.2495

$ws.Range("H131").Value = 2851.75
$ws.Range("I131").Value = 1473.4286
$ws.Range("K131").Value = 4420.2858
$ws.Range("M131").Value = 619.7142000000003

$ws.Range("H132").Value = 1284.3077
$ws.Range("I132").Value = 1220.3
$ws.Range("K132").Value = 3660.9
$ws.Range("M132").Value = -1130.9

$ws.Range("H135").Value = 658.5
$ws.Range("I135").Value = 565.2632
$ws.Range("K135").Value = 5087.3688
$ws.Range("M135").Value = -2552.3688

$ws.Range("H138").Value = 3616.1292
$ws.Range("I138").Value = 3130.4614
$ws.Range("K138").Value = 9391.3842
$ws.Range("M138").Value = -4251.3842

$ws.Range("H141").Value = 1366.125
$ws.Range("I141").Value = 1488.6666
$ws.Range("K141").Value = 4465.9998
$ws.Range("M141").Value = 714.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3808.9473
$ws.Range("I94").Value = 4898.7856
$ws.Range("K94").Value = 4898.7856
$ws.Range("M94").Value = -4447.7856

$ws.Range("H99").Value = 1878.35
$ws.Range("I99").Value = 1175.6923
$ws.Range("J99").Value = 3183.2856
$ws.Range("K99").Value = 1175.6923
$ws.Range("L99").Value = 3183.2856
$ws.Range("M99").Value = 322.3077000000001
$ws.Range("N99").Value = -6179.2856

$ws.Range("H107").Value = 2040.8
$ws.Range("I107").Value = 2126.5
$ws.Range("K107").Value = 2126.5
$ws.Range("M107").Value = -206.5

$ws.Range("H134").Value = 4304.75
$ws.Range("I134").Value = 4304.75
$ws.Range("K134").Value = 12914.25
$ws.Range("M134").Value = -10379.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 26805.3
$ws.Range("I22").Value = 2256.125
$ws.Range("J22").Value = 125002
$ws.Range("K22").Value = 2256.125
$ws.Range("L22").Value = 125002
$ws.Range("M22").Value = -1906.125
$ws.Range("N22").Value = -125702

$ws.Range("H31").Value = 2468.3333
$ws.Range("I31").Value = 1904.2858
$ws.Range("K31").Value = 1904.2858
$ws.Range("M31").Value = -1609.2858

$ws.Range("H34").Value = 2468.3333
$ws.Range("I34").Value = 1904.2858
$ws.Range("K34").Value = 1904.2858
$ws.Range("M34").Value = -1702.2858

$ws.Range("H58").Value = 2081.5386
$ws.Range("I58").Value = 1914.6364
$ws.Range("K58").Value = 1914.6364
$ws.Range("M58").Value = -1711.6364

$ws.Range("H62").Value = 2599.3333
$ws.Range("I62").Value = 2532.3333
$ws.Range("J62").Value = 2666.3333
$ws.Range("K62").Value = 2532.3333
$ws.Range("L62").Value = 2666.3333
$ws.Range("M62").Value = -1908.3333
$ws.Range("N62").Value = -3914.3333

$ws.Range("H65").Value = 2599.3333
$ws.Range("I65").Value = 2532.3333
$ws.Range("J65").Value = 2666.3333
$ws.Range("K65").Value = 12661.6665
$ws.Range("L65").Value = 13331.6665
$ws.Range("M65").Value = -9541.666499999999
$ws.Range("N65").Value = -19571.6665

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H134").Value = 3619
$ws.Range("I134").Value = 3465.5557
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 10396.6671
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -7861.667099999999
$ws.Range("N134").Value = -20070

$ws.Range("H136").Value = 2081.5386
$ws.Range("I136").Value = 1914.6364
$ws.Range("K136").Value = 5743.9092
$ws.Range("M136").Value = -3193.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3812.9167
$ws.Range("I132").Value = 5900
$ws.Range("J132").Value = 2322.1428
$ws.Range("K132").Value = 53100
$ws.Range("L132").Value = 20899.2852
$ws.Range("M132").Value = -50570
$ws.Range("N132").Value = -25959.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1143.4
$ws.Range("J102").Value = 3166.3333
$ws.Range("L102").Value = 3166.3333
$ws.Range("N102").Value = -6410.3333

$ws.Range("H107").Value = 4282.8335
$ws.Range("I107").Value = 849.5
$ws.Range("J107").Value = 5999.5
$ws.Range("K107").Value = 849.5
$ws.Range("L107").Value = 5999.5
$ws.Range("M107").Value = 1070.5
$ws.Range("N107").Value = -9839.5

$ws.Range("H122").Value = 4901.75
$ws.Range("I122").Value = 3199.6667
$ws.Range("J122").Value = 10008
$ws.Range("K122").Value = 9599.000100000001
$ws.Range("L122").Value = 30024
$ws.Range("M122").Value = -7149.000100000001
$ws.Range("N122").Value = -34924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2484.6667
$ws.Range("I93").Value = 2117.1
$ws.Range("K93").Value = 2117.1
$ws.Range("M93").Value = -869.0999999999999

$ws.Range("H100").Value = 3822.182
$ws.Range("I100").Value = 2409
$ws.Range("J100").Value = 4999.8335
$ws.Range("K100").Value = 2409
$ws.Range("L100").Value = 4999.8335
$ws.Range("M100").Value = -1868
$ws.Range("N100").Value = -6081.8335

$ws.Range("H101").Value = 46666.668
$ws.Range("J101").Value = 46666.668
$ws.Range("L101").Value = 46666.668
$ws.Range("N101").Value = -53156.668

$ws.Range("H132").Value = 4334.6665
$ws.Range("I132").Value = 3999.5
$ws.Range("K132").Value = 11998.5
$ws.Range("M132").Value = -9468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 21211
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 21211
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -22193

$ws.Range("H132").Value = 1462.625
$ws.Range("I132").Value = 1462.625
$ws.Range("K132").Value = 4387.875
$ws.Range("M132").Value = -1857.875

$ws.Range("H136").Value = 1280.2142
$ws.Range("I136").Value = 1313
$ws.Range("J136").Value = 395
$ws.Range("K136").Value = 3939
$ws.Range("L136").Value = 1185
$ws.Range("M136").Value = -1389
$ws.Range("N136").Value = -6285

